$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.918.90"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "'1.767.76"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'328.77"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.4530"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'0.3525"
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("D9").Value = "'41.98"
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("D10").Value = "'0.07382"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'1.094"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'20.71"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "'6.008"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "'7.181"
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("D16").Value = "'1.775.55"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").Value = "'92.58"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "'0.00001061"
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'0.06449"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'5.764"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'27.941.56"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'159.17"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").Value = "'20.14"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "'1.977.29"
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Value = "'2.149"
$ws.Range("E29").Value = "  +3.53%  "
$ws.Range("D30").Value = "'124.21"
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").Value = "'1.075"
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("D32").Value = "'0.09180"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("D33").Value = "'5.612"
$ws.Range("E33").Value = "  +2.45%  "
$ws.Range("D34").Value = "'3.664"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").Value = "'11.83"
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("D36").Value = "'0.02283"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.06115"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").Value = "'0.2092"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").Value = "'4.941"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "'0.6251"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "'1.177"
$ws.Range("E41").Value = "  -0.39%  "
$ws.Range("D42").Value = "'1.382"
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "'7.795"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").Value = "'13.21"
$ws.Range("E44").Value = "  +0.90%  "
$ws.Range("D45").Value = "'3.738"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").Value = "'0.5841"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "'122.39"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "'1.930"
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").Value = "'0.06829"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").Value = "'1.146"
$ws.Range("E51").Value = "  +1.96%  "

$ws.Range("D2:D51").Style = "Normal"
